$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 16 updates
# Plain filename text assigns cleanly without numeric coercion.
$ws.Range("D16").Value = "image_20250807111728_ppp0.jpg"

# The coordinate list looks numeric (commas) and would be auto-converted to a
# number by Excel, losing its comma formatting. Prefix with an apostrophe to
# force text entry, then reset the style so no extra/text number format is
# left behind on the cell (keeping it identical in style to the original).
$ws.Range("I16").Value = "'642,530,686,574"
$ws.Range("I16").Style = "Normal"

# Row 17 updates
$ws.Range("D17").Value = "image_20250807111728_ppp0.jpg"

$ws.Range("I17").Value = "'794,481,830,525"
$ws.Range("I17").Style = "Normal"

# "0.70" would be coerced to the number 0.7, dropping the trailing zero;
# force it to remain text the same way.
$ws.Range("J17").Value = "'0.70"
$ws.Range("J17").Style = "Normal"
